# Actualización automática 2026-01-01 08:30:07
# Insert two new client rows (alphabetically placed) into both report
# sheets ("VENTAS POR GRUPO" and "VENTA MENSUAL"), roll the monthly
# columns on "VENTA MENSUAL" forward by one month (sep-oct-nov-dic ->
# oct-nov-dic-ene), refresh the per-product/per-month figures, and widen
# the CLIENTE column to fit the longer new client names.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# ---------------------------------------------------------------------
# 1) Widen the CLIENTE column (B) on both sheets.
# ---------------------------------------------------------------------
$ws1.Columns.Item(2).ColumnWidth = 52
$ws2.Columns.Item(2).ColumnWidth = 52

# ---------------------------------------------------------------------
# 2) "VENTAS POR GRUPO" sheet
# ---------------------------------------------------------------------

# New client "DISTRIBUIDORA VASQUEZ ORDOÑEZ DISTRIVASOR CIA LTDA" is
# inserted alphabetically before "FABIMP BENIGNO BRAVO S.A.S." (row 10),
# pushing every following client row down by one.
$ws1.Rows.Item(10).Insert()
$ws1.Cells.Item(10, 1).Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Cells.Item(10, 2).Value = "DISTRIBUIDORA VASQUEZ ORDOÑEZ DISTRIVASOR CIA LTDA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(10, $col).Value = 0
}

# New client "VALENCIA RUIZ FAUSTO GABRIEL" is inserted alphabetically
# before "VIEJO RIVAS MAYRA ANABELLE" (now row 21 after the first insert).
$ws1.Rows.Item(21).Insert()
$ws1.Cells.Item(21, 1).Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Cells.Item(21, 2).Value = "VALENCIA RUIZ FAUSTO GABRIEL"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(21, $col).Value = 0
}

# Refresh this month's per-product figures: the clients that used to
# carry non-zero amounts sold nothing this period.
$ws1.Cells.Item(12, 5).Value = 0    # FRANK FERRETERIA FRANKFERRE CIA. - FREGADEROS DE COCINA
$ws1.Cells.Item(12, 9).Value = 0    # FRANK FERRETERIA FRANKFERRE CIA. - LAVABOS
$ws1.Cells.Item(12, 13).Value = 0   # FRANK FERRETERIA FRANKFERRE CIA. - PORCELANATO
$ws1.Cells.Item(14, 5).Value = 0    # ILLER LOPEZ ROBERTO FERNANDO - FREGADEROS DE COCINA
$ws1.Cells.Item(19, 12).Value = 0   # ROCAFUERTE LOPEZ EVELYN ESTEFANIA - PIEDRA SINTERIZADA
$ws1.Cells.Item(20, 5).Value = 0    # SARMIENTO SARMIENTO SANDRA EULALIA - FREGADEROS DE COCINA
$ws1.Cells.Item(20, 13).Value = 0   # SARMIENTO SARMIENTO SANDRA EULALIA - PORCELANATO

# Update the "X de 21" client-count summary row (now row 23, there are
# 21 clients, and none currently carries a non-zero figure).
$lastRow1 = 23
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item($lastRow1, $col).Value = "0 de 21"
}

$ws1.Range("A1:R" + $lastRow1).Select()

# ---------------------------------------------------------------------
# 3) "VENTA MENSUAL" sheet
# ---------------------------------------------------------------------

# Same two new client rows, in the same alphabetical slots.
$ws2.Rows.Item(10).Insert()
$ws2.Cells.Item(10, 1).Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Cells.Item(10, 2).Value = "DISTRIBUIDORA VASQUEZ ORDOÑEZ DISTRIVASOR CIA LTDA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(10, $col).Value = 0
}

$ws2.Rows.Item(21).Insert()
$ws2.Cells.Item(21, 1).Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Cells.Item(21, 2).Value = "VALENCIA RUIZ FAUSTO GABRIEL"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(21, $col).Value = 0
}

# Roll the monthly header labels forward by one month: the window now
# shows octubre / noviembre / diciembre / enero instead of septiembre /
# octubre / noviembre / diciembre.
$ws2.Cells.Item(1, 3).Value = "octubre"
$ws2.Cells.Item(1, 4).Value = "noviembre"
$ws2.Cells.Item(1, 5).Value = "diciembre"
$ws2.Cells.Item(1, 6).Value = "enero"

# Refresh each client's monthly figures for the new rolling window
# (columns C..F = octubre, noviembre, diciembre, enero).
$monthly = @{
    2  = @(0, 797.36, 1347.49, 0)       # ARMIJOS SALINAS LUIS CLAUDIO
    3  = @(0, 0, 0, 0)                  # ASES GAVILANEZ FAUSTO HERNAN
    4  = @(17.99, 0, 0, 0)              # BARROS YUNGA DIEGO VINICIO
    5  = @(0, 0, 0, 0)                  # BRAVO MONTENEGRO DANIEL ANDRES
    6  = @(0, 0, 0, 0)                  # BRITO CARDENAS RUTH CECILIA
    7  = @(0, 0, 0, 0)                  # COELLO TRONCOSO JOSE GREGORIO
    8  = @(0, 0, 0, 0)                  # COMERCIAL LUNA PAZMIÑO CIA. LTDA.
    9  = @(3992.9, 1631.15, 0, 0)       # CORPORACION AREVALO-YUMBLA E HIJOS
    10 = @(0, 0, 0, 0)                  # DISTRIBUIDORA VASQUEZ ORDOÑEZ DISTRIVASOR CIA LTDA
    11 = @(252.25, 0, 0, 0)             # FABIMP BENIGNO BRAVO S.A.S.
    12 = @(7662.57, 0, 4874.94, 0)      # FRANK FERRETERIA FRANKFERRE CIA.
    13 = @(0, 0, 6342.22, 0)            # HUERTA MUÑOZ NANCY ELIZABETH
    14 = @(58.48, 103.71, 93.98999999999999, 0)  # ILLER LOPEZ ROBERTO FERNANDO
    15 = @(3896.18, 0, 0, 0)            # MIM CONSTRUFERRETERIA E IMPORTADORA SAS
    16 = @(0, 0, 0, 0)                  # MOROCHO BACUILIMA HILDA INES
    17 = @(0, 0, 838.7, 0)              # MULLO GUACHO ANA LUCIA
    18 = @(326.73, 0, 0, 0)             # PAUTA ASTUDILLO JULIO HERNAN
    19 = @(1015.74, 1218.02, -591.61, 0) # ROCAFUERTE LOPEZ EVELYN ESTEFANIA
    20 = @(0, 2963.59, 0, 0)            # SARMIENTO SARMIENTO SANDRA EULALIA
    21 = @(0, 0, 232.74, 0)             # VALENCIA RUIZ FAUSTO GABRIEL
    22 = @(0, 0, 0, 0)                  # VIEJO RIVAS MAYRA ANABELLE
}

foreach ($row in $monthly.Keys) {
    $vals = $monthly[$row]
    $ws2.Cells.Item($row, 3).Value = $vals[0]
    $ws2.Cells.Item($row, 4).Value = $vals[1]
    $ws2.Cells.Item($row, 5).Value = $vals[2]
    $ws2.Cells.Item($row, 6).Value = $vals[3]
}

# Update the totals row (now row 23) for the new rolling window.
$lastRow2 = 23
$ws2.Cells.Item($lastRow2, 3).Value = 17222.84
$ws2.Cells.Item($lastRow2, 4).Value = 3750.24
$ws2.Cells.Item($lastRow2, 5).Value = 16102.06
$ws2.Cells.Item($lastRow2, 6).Value = 0
$ws2.Cells.Item($lastRow2, 7).Value = 0

$ws2.Range("A1:G" + $lastRow2).Select()

$ws1.Select()
